$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content so the sheet is rebuilt cleanly to the new shape
$ws.Cells.Clear()

$headers = @("Id", "Date", "Time", "League", "Home", "Away", "Odd_H_FT", "Odd_D_FT", "Odd_A_FT", "Odd_H_HT", "Odd_D_HT", "Odd_A_HT", "Odd_Over05_FT", "Odd_Under05_FT", "Odd_Over15_FT", "Odd_Under15_FT", "Odd_Over25_FT", "Odd_Under25_FT", "Odd_Over35_FT", "Odd_Under35_FT", "Odd_Over05_HT", "Odd_Under05_HT", "Odd_BTTS_Yes", "Odd_BTTS_No", "Odd_CS_1-0", "Odd_CS_2-0", "Odd_CS_2-1", "Odd_CS_3-0", "Odd_CS_3-1", "Odd_CS_3-2", "Odd_CS_0-0", "Odd_CS_1-1", "Odd_CS_2-2", "Odd_CS_3-3", "Odd_CS_4-4", "Odd_CS_0-1", "Odd_CS_0-2", "Odd_CS_1-2", "Odd_CS_0-3", "Odd_CS_1-3", "Odd_CS_2-3", "Odd_Over2_FT", "Odd_Under2_FT", "Odd_Over3_FT", "Odd_Under3_FT")
$row2 = @("2msHXImo", "25/03/2025", "16:45", "SCOTLAND - CHAMPIONSHIP", "Livingston", "Falkirk", 3, 3.2, 2.38, 3.6, 2, 3, 1.06, 8, 1.36, 3, 2.2, 1.65, 4, 1.22, 1.5, 2.5, 1.91, 1.8, 8.5, 15, 12, 29, 26, 41, 8, 6, 17, 51, 800, 7, 11, 10, 21, 21, 34, $null, $null, $null, $null)
$row3 = @("U1DzKp1e", "25/03/2025", "16:45", "SCOTLAND - CHAMPIONSHIP", "Morton", "Raith", 2.7, 2.8, 2.9, 3.5, 1.83, 3.6, 1.13, 5.5, 1.57, 2.25, 2.75, 1.4, 5.5, 1.13, 1.62, 2.2, 2.25, 1.57, 6.5, 11, 12, 26, 29, 41, 5.5, 6, 21, 81, 1000, 6.5, 12, 13, 29, 29, 41, 1.99, 1.74, $null, $null)
$row4 = @("QJJhSgNr", "25/03/2025", "20:00", "URUGUAY - LIGA AUF URUGUAYA", "Cerro Largo", "Montevideo City", 2.6, 3, 2.8, 3.2, 2.1, 3.4, 1.06, 10, 1.29, 3.5, 2.03, 1.83, 3.4, 1.3, 1.4, 2.75, 1.73, 2, 9, 13, 10, 26, 21, 29, 9.5, 6, 13, 41, 201, 9.5, 13, 11, 29, 23, 29, $null, $null, 1.88, 1.98)

$data = @($headers, $row2, $row3, $row4)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# Re-apply the header style (bold, bordered, centered/top-aligned) to the
# full header row, including the newly-added columns AP:AS
$headerRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1, $headers.Length))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$ws.Range("A1").Select()
